$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the old row 216, pushing old rows
# 216, 217, 218 down to 218, 219, 220 (their contents/styles travel with them).
$ws.Rows("216:217").Insert()

# --- Row 214: same record, but re-dated and re-volumed ---
$ws.Range("D214").Value = 44568
$ws.Range("J214").Value = 2000

# --- Row 215: re-dated, re-volumed, quality downgraded to Segunda, prices lowered ---
$ws.Range("D215").Value = 44568
$ws.Range("I215").Value = "Segunda"
$ws.Range("J215").Value = 2000
$ws.Range("K215").Value = 300
$ws.Range("L215").Value = 300
$ws.Range("M215").Value = 300
$ws.Range("P215").Value = 300

# --- Row 216: newly inserted row, carries what used to be row 214's data ---
$ws.Range("A216").Value = 5
$ws.Range("B216").Value = "Macroferia Regional de Talca"
$ws.Range("C216").Value = "Maule"
$ws.Range("D216").Value = 44357
$ws.Range("E216").Value = 7
$ws.Range("F216").Value = 100112006
$ws.Range("G216").Value = "Repollo"
$ws.Range("H216").Value = "Crespo record"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 5000
$ws.Range("K216").Value = 500
$ws.Range("L216").Value = 500
$ws.Range("M216").Value = 500
$ws.Range("N216").Value = "$/unidad"
$ws.Range("O216").Value = "Región del Maule"
$ws.Range("P216").Value = 500
$ws.Range("Q216").Value = 1
$ws.Range("R216").Value = "Hortaliza"

# --- Row 217: newly inserted row, carries what used to be row 215's data ---
$ws.Range("A217").Value = 5
$ws.Range("B217").Value = "Macroferia Regional de Talca"
$ws.Range("C217").Value = "Maule"
$ws.Range("D217").Value = 44371
$ws.Range("E217").Value = 7
$ws.Range("F217").Value = 100112006
$ws.Range("G217").Value = "Repollo"
$ws.Range("H217").Value = "Crespo record"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 5000
$ws.Range("K217").Value = 450
$ws.Range("L217").Value = 450
$ws.Range("M217").Value = 450
$ws.Range("N217").Value = "$/unidad"
$ws.Range("O217").Value = "Región del Maule"
$ws.Range("P217").Value = 450
$ws.Range("Q217").Value = 1
$ws.Range("R217").Value = "Hortaliza"
